# Swap the data of row 8 and row 9 (everything except the few columns
# that already hold identical values in both rows: C, D, J, K, N, S, T,
# U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY).
#
# Columns that actually differ between row 8 and row 9 and therefore
# need to be swapped: A, B, E, F, G, H, I, P, Q, R.
#
# Column I ("Antal") is stored as text (e.g. "2", "7") even though the
# text looks numeric, so a plain Value2 assignment of a numeric-looking
# string would get auto-coerced to a real number by Excel. Using
# Copy + PasteSpecial(xlPasteValues) moves the cell's value (and its
# underlying type) verbatim, without forcing any number-format / style
# change and without the auto type coercion that a direct Value2
# assignment of a numeric-looking string would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$cols = @("A", "B", "E", "F", "G", "H", "I", "P", "Q", "R")

foreach ($col in $cols) {
    $addrRow8 = "${col}8"
    $addrRow9 = "${col}9"
    $addrTemp = "${col}1000"

    # stash row 8's value
    $ws.Range($addrRow8).Copy() | Out-Null
    $ws.Range($addrTemp).PasteSpecial($xlPasteValues) | Out-Null

    # row 9 -> row 8
    $ws.Range($addrRow9).Copy() | Out-Null
    $ws.Range($addrRow8).PasteSpecial($xlPasteValues) | Out-Null

    # stashed row 8 -> row 9
    $ws.Range($addrTemp).Copy() | Out-Null
    $ws.Range($addrRow9).PasteSpecial($xlPasteValues) | Out-Null

    # clean up the temp holding cell
    $ws.Range($addrTemp).ClearContents() | Out-Null
}

$excel.CutCopyMode = 0
